# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (fund-level holdings detail) right
#    before the "总计" (totals) summary sheet - cloned from "2021-Q4" so it
#    inherits the exact same column layout/styles, then its content is
#    overwritten with the 2022-Q1 figures.
# 2) Update the "总计" sheet: prepend a new 2022-Q1 summary row and shift
#    the previously existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: build the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetBeforeCopy = $wb.Worksheets.Item("总计")

# Clone "2021-Q4" so the new sheet lands right before "总计" with matching
# column widths / header & border styling, then rename it.
$srcSheet.Copy($totalSheetBeforeCopy)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# NOTE: inserting/copying a sheet invalidates previously-fetched worksheet
# handles in this runtime, so re-resolve "总计" by name now that the new
# sheet has been inserted in front of it.
$totalSheet = $wb.Worksheets.Item("总计")

# "2021-Q4" has 12 data rows (rows 2-13); 2022-Q1 only needs 11 (rows 2-12).
$q1.Rows.Item(13).Delete()

# Columns whose values look numeric but must stay TEXT (leading zeros in
# fund codes, percentage-like strings) - pre-format as Text so the literal
# strings are not silently coerced into numbers.
$q1.Range("B2:B12").NumberFormat = "@"
$q1.Range("D2:G12").NumberFormat = "@"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "010717"
$q1.Range("C2").Value = "前海开源优质企业6个月持有期混合A"
$q1.Range("D2").Value = "58.52"
$q1.Range("E2").Value = "93.93"
$q1.Range("F2").Value = "3.67"
$q1.Range("G2").Value = "2.1477"
$q1.Range("H2").Value = 10

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "008188"
$q1.Range("C3").Value = "前海开源稳健增长三年持有期混合"
$q1.Range("D3").Value = "25.26"
$q1.Range("E3").Value = "67.84"
$q1.Range("F3").Value = "4.48"
$q1.Range("G3").Value = "1.1316"
$q1.Range("H3").Value = 6

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "001837"
$q1.Range("C4").Value = "前海开源沪港深蓝筹精选灵活配置混合"
$q1.Range("D4").Value = "15.15"
$q1.Range("E4").Value = "94.10"
$q1.Range("F4").Value = "5.00"
$q1.Range("G4").Value = "0.7575"
$q1.Range("H4").Value = 10

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "010718"
$q1.Range("C5").Value = "前海开源优质企业6个月持有期混合C"
$q1.Range("D5").Value = "7.00"
$q1.Range("E5").Value = "93.93"
$q1.Range("F5").Value = "3.67"
$q1.Range("G5").Value = "0.2569"
$q1.Range("H5").Value = 10

$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "006775"
$q1.Range("C6").Value = "前海开源优质成长混合"
$q1.Range("D6").Value = "3.68"
$q1.Range("E6").Value = "68.45"
$q1.Range("F6").Value = "4.43"
$q1.Range("G6").Value = "0.1630"
$q1.Range("H6").Value = 6

$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "011287"
$q1.Range("C7").Value = "前海开源聚慧三年持有期混合"
$q1.Range("D7").Value = "3.25"
$q1.Range("E7").Value = "68.77"
$q1.Range("F7").Value = "4.31"
$q1.Range("G7").Value = "0.1401"
$q1.Range("H7").Value = 7

$q1.Range("A8").Value = 6
$q1.Range("B8").Value = "006216"
$q1.Range("C8").Value = "前海开源价值成长灵活配置混合A"
$q1.Range("D8").Value = "1.97"
$q1.Range("E8").Value = "67.45"
$q1.Range("F8").Value = "4.23"
$q1.Range("G8").Value = "0.0833"
$q1.Range("H8").Value = 7

$q1.Range("A9").Value = 7
$q1.Range("B9").Value = "011018"
$q1.Range("C9").Value = "景顺长城安泽回报一年持有期混合A"
$q1.Range("D9").Value = "7.81"
$q1.Range("E9").Value = "34.82"
$q1.Range("F9").Value = "0.87"
$q1.Range("G9").Value = "0.0679"
$q1.Range("H9").Value = 9

$q1.Range("A10").Value = 8
$q1.Range("B10").Value = "006217"
$q1.Range("C10").Value = "前海开源价值成长灵活配置混合C"
$q1.Range("D10").Value = "0.73"
$q1.Range("E10").Value = "67.45"
$q1.Range("F10").Value = "4.23"
$q1.Range("G10").Value = "0.0309"
$q1.Range("H10").Value = 7

$q1.Range("A11").Value = 9
$q1.Range("B11").Value = "011019"
$q1.Range("C11").Value = "景顺长城安泽回报一年持有期混合C"
$q1.Range("D11").Value = "0.24"
$q1.Range("E11").Value = "34.82"
$q1.Range("F11").Value = "0.87"
$q1.Range("G11").Value = "0.0021"
$q1.Range("H11").Value = 9

$q1.Range("A12").Value = 10
$q1.Range("B12").Value = "002860"
$q1.Range("C12").Value = "前海开源沪港深新机遇灵活配置混合"
$q1.Range("D12").Value = "0.01"
$q1.Range("E12").Value = "83.26"
$q1.Range("F12").Value = "6.14"
$q1.Range("G12").Value = "0.0006"
$q1.Range("H12").Value = 10

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet - add the 2022-Q1 row, shift the rest
# ---------------------------------------------------------------------
# Extend the styled "index" column (A) down into the new row 7 first by
# cloning the format of the last existing row, so it keeps the bold/
# bordered look of A2:A6.
$totalSheet.Range("A6").Copy($totalSheet.Range("A7"))

# Rewrite rows bottom-up isn't required since every cell is assigned its
# final literal value directly (no relative/formula shifting involved).
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 11
$totalSheet.Range("D2").Value = 4.78

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 12
$totalSheet.Range("D3").Value = 4.3

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 12
$totalSheet.Range("D4").Value = 9.43

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 3
$totalSheet.Range("D5").Value = 0.23

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q1"
$totalSheet.Range("C6").Value = 2
$totalSheet.Range("D6").Value = 0.19

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2020-Q4"
$totalSheet.Range("C7").Value = 2
$totalSheet.Range("D7").Value = 0.26

# Restore the originally-active first tab (inserting/renaming sheets above
# shifted the active-sheet selection onto the new "2022-Q1" sheet).
$wb.Worksheets.Item(1).Activate()
